$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column A for rows 2-4
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 7
$ws.Range("A4").Value = 20

# Change selection to A5 (reflected as sheetView selection activeCell/sqref)
$ws.Range("A5").Select()
